$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Alfresco folder path text: "CaseFiles" -> "Case Files"
$ws.Range("D16").Value = "setEcmFolderPath, '/Sites/acm/documentLibrary/Case Files/' + dateFormat('yyyyMMdd') + '_' + `$caseFile.getId()"

# Add a new "Set Case Status" rule row
$ws.Range("B17").Value = "Set Case Status"
$ws.Range("C17").Value = "status == null"
$ws.Range("D17").Value = "setStatus, 'DRAFT'"

# Update the current selection to D16
$ws.Range("D16").Select()
